$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 185
$ws.Range("J2").Value = 763
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 204
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 129
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 86
$ws.Range("T2").Value = 147
$ws.Range("V2").Value = 1126
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1149
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 18
$ws.Range("AA2").Value = 6
